$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat")

# --- Title / header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 32   Number  13"
$ws.Range("C9").Value = "Report Covering the Week  3/24/2025  Through  3/30/2025"

# --- CompStat data table updates ---
$ws.Range("C14").Value = 1
$ws.Range("C14").NumberFormat = '#,##0'
$ws.Range("F14").Value = 1
$ws.Range("F14").NumberFormat = '#,##0'
$ws.Range("I14").Value = 2
$ws.Range("K14").Value = 100
$ws.Range("L14").Value = -60
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = -66.666666666666
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 7
$ws.Range("G15").Value = 6
$ws.Range("H15").Value = 16.666666666666
$ws.Range("I15").Value = 16
$ws.Range("J15").Value = 15
$ws.Range("K15").Value = 6.666666666666
$ws.Range("L15").Value = 33.333333333333
$ws.Range("M15").Value = 45.454545454545
$ws.Range("N15").Value = -36
$ws.Range("D16").Value = 15
$ws.Range("E16").Value = -80
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 24
$ws.Range("H16").Value = -50
$ws.Range("I16").Value = 62
$ws.Range("J16").Value = 73
$ws.Range("K16").Value = -15.068493150684
$ws.Range("L16").Value = -28.735632183908
$ws.Range("M16").Value = -27.906976744186
$ws.Range("N16").Value = -82.183908045977
$ws.Range("C17").Value = 16
$ws.Range("D17").Value = 19
$ws.Range("E17").Value = -15.789473684210
$ws.Range("F17").Value = 73
$ws.Range("G17").Value = 94
$ws.Range("H17").Value = -22.340425531914
$ws.Range("I17").Value = 209
$ws.Range("J17").Value = 254
$ws.Range("K17").Value = -17.716535433070
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 68.548387096774
$ws.Range("N17").Value = -13.278008298755
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -57.142857142857
$ws.Range("F18").Value = 18
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = -5.263157894736
$ws.Range("I18").Value = 102
$ws.Range("J18").Value = 64
$ws.Range("K18").Value = 59.375
$ws.Range("L18").Value = 6.25
$ws.Range("M18").Value = -15.702479338843
$ws.Range("N18").Value = -87.218045112782
$ws.Range("C19").Value = 24
$ws.Range("D19").Value = 33
$ws.Range("E19").Value = -27.272727272727
$ws.Range("F19").Value = 98
$ws.Range("G19").Value = 123
$ws.Range("H19").Value = -20.325203252032
$ws.Range("I19").Value = 296
$ws.Range("J19").Value = 395
$ws.Range("K19").Value = -25.063291139240
$ws.Range("L19").Value = -19.565217391304
$ws.Range("M19").Value = 41.626794258373
$ws.Range("N19").Value = -17.777777777777
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -80
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 19
$ws.Range("H20").Value = -63.157894736842
$ws.Range("J20").Value = 56
$ws.Range("K20").Value = -21.428571428571
$ws.Range("L20").Value = -50
$ws.Range("M20").Value = -36.231884057971
$ws.Range("N20").Value = -96.372629843363
$ws.Range("C21").Value = 49
$ws.Range("D21").Value = 80
$ws.Range("E21").Value = -38.75
$ws.Range("F21").Value = 216
$ws.Range("G21").Value = 285
$ws.Range("H21").Value = -24.210526315789
$ws.Range("I21").Value = 731
$ws.Range("J21").Value = 858
$ws.Range("K21").Value = -14.801864801864
$ws.Range("L21").Value = -15.491329479768
$ws.Range("M21").Value = 17.524115755627
$ws.Range("N21").Value = -75.560013373453
$ws.Range("C23").Value = 1
$ws.Range("E23").Value = -83.333333333333
$ws.Range("F23").Value = 10
$ws.Range("G23").Value = 14
$ws.Range("H23").Value = -28.571428571428
$ws.Range("I23").Value = 28
$ws.Range("J23").Value = 24
$ws.Range("K23").Value = 16.666666666666
$ws.Range("L23").Value = -34.883720930232
$ws.Range("M23").Value = 154.545454545455
$ws.Range("C24").Value = 89
$ws.Range("D24").Value = 93
$ws.Range("E24").Value = -4.301075268817
$ws.Range("F24").Value = 305
$ws.Range("G24").Value = 363
$ws.Range("H24").Value = -15.977961432506
$ws.Range("I24").Value = 1058
$ws.Range("J24").Value = 1104
$ws.Range("K24").Value = -4.166666666666
$ws.Range("L24").Value = 8.959835221421
$ws.Range("M24").Value = 35.294117647058
$ws.Range("C25").Value = 55
$ws.Range("D25").Value = 62
$ws.Range("E25").Value = -11.290322580645
$ws.Range("F25").Value = 177
$ws.Range("G25").Value = 221
$ws.Range("H25").Value = -19.909502262443
$ws.Range("I25").Value = 631
$ws.Range("J25").Value = 630
$ws.Range("K25").Value = 0.158730158730
$ws.Range("L25").Value = 32.563025210084
$ws.Range("C26").Value = 44
$ws.Range("D26").Value = 34
$ws.Range("E26").Value = 29.411764705882
$ws.Range("F26").Value = 155
$ws.Range("G26").Value = 134
$ws.Range("H26").Value = 15.671641791044
$ws.Range("I26").Value = 433
$ws.Range("J26").Value = 408
$ws.Range("K26").Value = 6.127450980392
$ws.Range("L26").Value = 6.913580246913
$ws.Range("M26").Value = -0.230414746543
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = 57.142857142857
$ws.Range("I27").Value = 24
$ws.Range("J27").Value = 27
$ws.Range("K27").Value = -11.111111111111
$ws.Range("L27").Value = 33.333333333333
$ws.Range("C28").Value = 4
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 12
$ws.Range("H28").Value = -7.692307692307
$ws.Range("I28").Value = 52
$ws.Range("J28").Value = 42
$ws.Range("K28").Value = 23.809523809523
$ws.Range("L28").Value = 10.638297872340
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 100
$ws.Range("I29").Value = 3
$ws.Range("K29").Value = -40
$ws.Range("L29").Value = -25
$ws.Range("M29").Value = -66.666666666666
$ws.Range("N29").Value = -85
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 100
$ws.Range("I30").Value = 3
$ws.Range("K30").Value = -40
$ws.Range("L30").Value = -25
$ws.Range("M30").Value = -62.5
$ws.Range("N30").Value = -83.333333333333
$ws.Range("D31").Value = 1
$ws.Range("D31").NumberFormat = '#,##0'
$ws.Range("E31").Value = -100
$ws.Range("E31").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F31").Value = 1
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = -50
$ws.Range("J31").Value = 7
$ws.Range("K31").Value = -71.428571428571
$ws.Range("L33").Value = -50
